$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/14/2025  Through  7/20/2025"

# --- Column E width update (match column H bestfit width) ---
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# --- Cells changing between numeric and text representation ---

# G14: numeric 1 -> text "0" (reuses style of C14 / shared string "0")
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G14").Formula = '=TEXT(0,"0")'
$ws.Range("G14").Copy()
$ws.Range("G14").PasteSpecial(-4163)

# H14: numeric -100 -> text "***.*" (reuses style of E14 / shared string "***.*")
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("H14").Value = "***.*"

# C28: numeric 1 -> text "0" (reuses style of D29 / shared string "0")
$ws.Range("D29").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Formula = '=TEXT(0,"0")'
$ws.Range("C28").Copy()
$ws.Range("C28").PasteSpecial(-4163)

# C29: text "0" -> numeric 1 (reuses style of D28)
$ws.Range("D28").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1

# F29: text "0" -> numeric 1 (reuses style of D28)
$ws.Range("D28").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value = 1

# C30: text "0" -> numeric 1 (reuses style of D28)
$ws.Range("D28").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 1

# F30: text "0" -> numeric 1 (reuses style of D28)
$ws.Range("D28").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").Value = 1

# --- Simple same-type/style numeric value updates ---
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("N15").Value = 29.411764705882
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 6
$ws.Range("F16").Value = 25
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 138
$ws.Range("J16").Value = 117
$ws.Range("K16").Value = 17.948717948717
$ws.Range("L16").Value = 112.307692307692
$ws.Range("M16").Value = 126.229508196721
$ws.Range("N16").Value = -81.043956043956
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -4
$ws.Range("I17").Value = 190
$ws.Range("J17").Value = 157
$ws.Range("K17").Value = 21.019108280254
$ws.Range("L17").Value = 84.466019417475
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = -32.624113475177
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 108
$ws.Range("J18").Value = 111
$ws.Range("K18").Value = -2.702702702702
$ws.Range("L18").Value = 47.945205479452
$ws.Range("M18").Value = 2.857142857142
$ws.Range("N18").Value = -90.90909090909
$ws.Range("C19").Value = 33
$ws.Range("D19").Value = 36
$ws.Range("E19").Value = -8.333333333333
$ws.Range("F19").Value = 139
$ws.Range("G19").Value = 150
$ws.Range("H19").Value = -7.333333333333
$ws.Range("I19").Value = 952
$ws.Range("J19").Value = 986
$ws.Range("K19").Value = -3.448275862068
$ws.Range("L19").Value = -5.084745762711
$ws.Range("M19").Value = 10.826542491268
$ws.Range("N19").Value = -75.941369724538
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 38
$ws.Range("J20").Value = 28
$ws.Range("K20").Value = 35.714285714285
$ws.Range("L20").Value = -5
$ws.Range("M20").Value = 58.333333333333
$ws.Range("N20").Value = -85.66037735849
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 49
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 209
$ws.Range("G21").Value = 223
$ws.Range("H21").Value = -6.278026905829
$ws.Range("I21").Value = 1448
$ws.Range("J21").Value = 1406
$ws.Range("K21").Value = 2.987197724039
$ws.Range("L21").Value = 12.074303405572
$ws.Range("M21").Value = 27.352682497801
$ws.Range("N21").Value = -77.536456717344
$ws.Range("C22").Value = 1
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 41
$ws.Range("J22").Value = 51
$ws.Range("K22").Value = -19.607843137254
$ws.Range("L22").Value = -4.651162790697
$ws.Range("M22").Value = 41.379310344827
$ws.Range("C24").Value = 56
$ws.Range("D24").Value = 88
$ws.Range("E24").Value = -36.363636363636
$ws.Range("G24").Value = 273
$ws.Range("H24").Value = -27.472527472527
$ws.Range("I24").Value = 1395
$ws.Range("J24").Value = 1664
$ws.Range("K24").Value = -16.165865384615
$ws.Range("L24").Value = -4.843110504774
$ws.Range("M24").Value = 38.118811881188
$ws.Range("C25").Value = 46
$ws.Range("D25").Value = 82
$ws.Range("E25").Value = -43.90243902439
$ws.Range("F25").Value = 170
$ws.Range("G25").Value = 271
$ws.Range("H25").Value = -37.269372693726
$ws.Range("I25").Value = 1286
$ws.Range("J25").Value = 1612
$ws.Range("K25").Value = -20.223325062034
$ws.Range("L25").Value = -13.108108108108
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 59
$ws.Range("G26").Value = 66
$ws.Range("H26").Value = -10.60606060606
$ws.Range("I26").Value = 403
$ws.Range("J26").Value = 384
$ws.Range("K26").Value = 4.947916666666
$ws.Range("L26").Value = 2.544529262086
$ws.Range("M26").Value = 56.8093385214
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = -60
$ws.Range("J28").Value = 51
$ws.Range("K28").Value = 25.490196078431
$ws.Range("L28").Value = 20.754716981132
$ws.Range("I29").Value = 2
$ws.Range("K29").Value = -33.333333333333
$ws.Range("L29").Value = -50
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = -66.666666666666
$ws.Range("I30").Value = 2
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = -33.333333333333
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = -50
